$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Range("C19").Value = 92
